$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, [string]$val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$ws.Range("D2").Value = "61.514.11"
$ws.Range("E2").Value = "  +1.20%  "
$ws.Range("D3").Value = "3.445.40"
$ws.Range("E3").Value = "  +1.87%  "
Set-TextValue "D4" "1.00"
$ws.Range("E4").Value = "  -0.09%  "
Set-TextValue "D5" "580.05"
$ws.Range("E5").Value = "  +0.97%  "
Set-TextValue "D6" "148.59"
$ws.Range("E6").Value = "  +8.74%  "
$ws.Range("D7").Value = "3.447.89"
$ws.Range("E7").Value = "  +1.97%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  +0.79%  "
$ws.Range("E10").Value = "  +2.94%  "
$ws.Range("E11").Value = "  +1.50%  "
$ws.Range("E12").Value = "  +0.68%  "
$ws.Range("D13").Value = "4.035.38"
$ws.Range("E13").Value = "  +1.73%  "
Set-TextValue "D14" "28.03"
$ws.Range("E14").Value = "  +6.27%  "
$ws.Range("E15").Value = "  -0.41%  "
$ws.Range("E16").Value = "  +1.36%  "
$ws.Range("D17").Value = "3.439.64"
$ws.Range("E17").Value = "  +1.37%  "
$ws.Range("D18").Value = "61.604.91"
$ws.Range("E18").Value = "  +1.09%  "
Set-TextValue "D19" "6.32"
$ws.Range("E19").Value = "  +8.28%  "
$ws.Range("E20").Value = "  +2.66%  "
Set-TextValue "D21" "9.45"
$ws.Range("E21").Value = "  +0.11%  "
Set-TextValue "D22" "387.59"
$ws.Range("E22").Value = "  +2.72%  "
$ws.Range("E23").Value = "  +2.54%  "
$ws.Range("D24").Value = "3.592.17"
$ws.Range("E24").Value = "  +1.79%  "
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue "D25" "1.00"
$ws.Range("E25").Value = "  +0.32%  "
$ws.Range("B26").Value = "Litecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue "D26" "72.81"
$ws.Range("E26").Value = "  +1.98%  "
Set-TextValue "D27" "5.77"
$ws.Range("E27").Value = "  +0.68%  "
$ws.Range("E28").Value = "  -1.80%  "
$ws.Range("E29").Value = "  +6.81%  "
$ws.Range("E30").Value = "  +3.29%  "
Set-TextValue "D31" "0.999"
$ws.Range("E31").Value = "  -0.10%  "
Set-TextValue "D32" "1.53"
$ws.Range("E32").Value = "  -14.76%  "
Set-TextValue "D33" "8.26"
$ws.Range("E33").Value = "  +0.87%  "
$ws.Range("E34").Value = "  +0.75%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("E36").Value = "  +0.91%  "
$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D37" "5.27"
$ws.Range("E37").Value = "  +1.16%  "
$ws.Range("B38").Value = "Aptos"
$ws.Range("C38").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D38" "7.08"
$ws.Range("E38").Value = "  +3.05%  "
$ws.Range("E39").Value = "  +1.93%  "
Set-TextValue "D40" "166.18"
$ws.Range("E40").Value = "  +0.97%  "
$ws.Range("E41").Value = "  +4.75%  "
Set-TextValue "D42" "26.22"
$ws.Range("E42").Value = "  +7.16%  "
$ws.Range("E43").Value = "  +2.87%  "
$ws.Range("E44").Value = "  -0.09%  "
$ws.Range("E45").Value = "  +1.91%  "
Set-TextValue "D46" "42.29"
$ws.Range("E46").Value = "  +1.84%  "
Set-TextValue "D47" "1.73"
$ws.Range("E47").Value = "  +0.83%  "
$ws.Range("D48").Value = "2.613.76"
$ws.Range("E48").Value = "  +8.48%  "
Set-TextValue "D49" "1.17"
$ws.Range("E49").Value = "  -2.90%  "
$ws.Range("E50").Value = "  +3.39%  "
Set-TextValue "D51" "23.12"
$ws.Range("E51").Value = "  -1.39%  "

Write-Output "Applied 87 cell updates"
